$d = $word.ActiveDocument

$replacements = @(
    @("93×25=", "17×54="),
    @("34×92=", "91×92="),
    @("61×27=", "33×56="),
    @("62×85=", "29×68="),
    @("57×47=", "16×39="),
    @("94×64=", "46×37="),
    @("33×35=", "88×90="),
    @("37×92=", "48×63="),
    @("82×94=", "81×61="),
    @("47×52=", "78×78="),
    @("87×76=", "96×44="),
    @("78×93=", "50×99="),
    @("49×40=", "20×89="),
    @("49×84=", "67×84="),
    @("62×52=", "81×55="),
    @("81×43=", "39×96="),
    @("49×17=", "59×27="),
    @("92×27=", "61×17="),
    @("19×51=", "22×89="),
    @("54×35=", "54×37="),
    @("49×43=", "17×72="),
    @("40×66=", "91×63="),
    @("38×88=", "36×64="),
    @("50×36=", "16×40="),
    @("68×28=", "40×94=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
